$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.458.02"
$ws.Range("E2").Value = "  +1.85%  "

$ws.Range("D3").Value = "1.863.90"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3802"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07329"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9373"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07806"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("D13").Value = "1.857.79"
$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.449"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.560"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008902"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.04%  "

$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "27.471.52"
$ws.Range("E20").Value = "  +1.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.117"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.17%  "

$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.940"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("E27").Value = "  +1.36%  "

$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.960"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08896"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.216"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7592"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.605"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.729"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02061"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5629"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05287"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.990"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.074"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.696"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1531"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4906"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.659"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("E49").Value = "  +2.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06087"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9212"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.80%  "
